$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 'Price' column (D) holds values that read like numbers (e.g. '0.550',
# '18.30') but must stay literal TEXT -- exactly as the source data has them --
# including any trailing zeros that a numeric cast would silently drop.
# Forcing NumberFormat to text ('@') before the write keeps the digits intact,
# then resetting Style back to 'Normal' collapses the cell back onto the
# workbook's default style so no stray formatting is introduced.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '42.928.43'
$ws.Range("E2").Value = '  -1.82%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.572.71'
$ws.Range("E3").Value = '  +1.44%  '

# Row 4
Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.17%  '

# Row 5
Set-TextValue $ws.Range("D5") '302.35'
$ws.Range("E5").Value = '  +0.77%  '

# Row 6
Set-TextValue $ws.Range("D6") '97.66'
$ws.Range("E6").Value = '  +2.28%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.576'
$ws.Range("E7").Value = '  -0.31%  '

# Row 8
$ws.Range("E8").Value = '  -0.16%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.550'
$ws.Range("E9").Value = '  -1.33%  '

# Row 10
Set-TextValue $ws.Range("D10") '36.39'
$ws.Range("E10").Value = '  -1.19%  '

# Row 11
$ws.Range("E11").Value = '  -0.04%  '

# Row 12
$ws.Range("E12").Value = '  -0.40%  '

# Row 13
$ws.Range("E13").Value = '  +6.25%  '

# Row 14
Set-TextValue $ws.Range("D14") '2.577.28'
$ws.Range("E14").Value = '  +0.28%  '

# Row 15
Set-TextValue $ws.Range("D15") '0.887'
$ws.Range("E15").Value = '  +0.57%  '

# Row 16
Set-TextValue $ws.Range("D16") '14.41'
$ws.Range("E16").Value = '  +1.29%  '

# Row 17
Set-TextValue $ws.Range("D17") '42.941.76'
$ws.Range("E17").Value = '  -1.75%  '

# Row 18
Set-TextValue $ws.Range("D18") '0.0₃0998'
$ws.Range("E18").Value = '  +2.17%  '

# Row 19
Set-TextValue $ws.Range("D19") '12.92'
$ws.Range("E19").Value = '  +3.98%  '

# Row 20
Set-TextValue $ws.Range("D20") '6.64'
$ws.Range("E20").Value = '  +0.13%  '

# Row 21
Set-TextValue $ws.Range("D21") '72.09'
$ws.Range("E21").Value = '  -1.45%  '

# Row 22
Set-TextValue $ws.Range("D22") '254.87'
$ws.Range("E22").Value = '  -3.19%  '

# Row 23
Set-TextValue $ws.Range("D23") '2.96'
$ws.Range("E23").Value = '  +1.53%  '

# Row 24
Set-TextValue $ws.Range("D24") '2.13'
$ws.Range("E24").Value = '  -2.66%  '

# Row 25
Set-TextValue $ws.Range("D25") '28.81'
$ws.Range("E25").Value = '  -0.37%  '

# Row 26
$ws.Range("E26").Value = '  +0.02%  '

# Row 27
Set-TextValue $ws.Range("D27") '10.28'
$ws.Range("E27").Value = '  +0.99%  '

# Row 28
Set-TextValue $ws.Range("D28") '37.96'
$ws.Range("E28").Value = '  +0.63%  '

# Row 29
Set-TextValue $ws.Range("D29") '2.12'
$ws.Range("E29").Value = '  -5.14%  '

# Row 30
Set-TextValue $ws.Range("D30") '6.05'
$ws.Range("E30").Value = '  -0.89%  '

# Row 31
Set-TextValue $ws.Range("D31") '155.21'
$ws.Range("E31").Value = '  +2.29%  '

# Row 32
Set-TextValue $ws.Range("D32") '3.41'
$ws.Range("E32").Value = '  -3.20%  '

# Row 33
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D33") '2.16'
$ws.Range("E33").Value = '  +0.00%  '

# Row 34
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D34") '2.75'
$ws.Range("E34").Value = '  -1.11%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.0805'
$ws.Range("E35").Value = '  -0.24%  '

# Row 36
Set-TextValue $ws.Range("D36") '18.30'
$ws.Range("E36").Value = '  +10.89%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.114'
$ws.Range("E37").Value = '  -1.88%  '

# Row 38
$ws.Range("E38").Value = '  +0.04%  '

# Row 39
Set-TextValue $ws.Range("D39") '23.18'
$ws.Range("E39").Value = '  -1.73%  '

# Row 40
Set-TextValue $ws.Range("D40") '3.44'
$ws.Range("E40").Value = '  -2.99%  '

# Row 41
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D41") '3.89'
$ws.Range("E41").Value = '  +1.54%  '

# Row 42
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D42") '0.0311'
$ws.Range("E42").Value = '  -0.84%  '

# Row 43
Set-TextValue $ws.Range("D43") '2.07'
$ws.Range("E43").Value = '  +26.46%  '

# Row 44
Set-TextValue $ws.Range("D44") '2.073.13'
$ws.Range("E44").Value = '  +2.68%  '

# Row 45
Set-TextValue $ws.Range("D45") '0.998'
$ws.Range("E45").Value = '  -0.04%  '

# Row 46
Set-TextValue $ws.Range("D46") '9.21'
$ws.Range("E46").Value = '  +1.69%  '

# Row 47
Set-TextValue $ws.Range("D47") '85.47'
$ws.Range("E47").Value = '  -1.98%  '

# Row 48
Set-TextValue $ws.Range("D48") '76.38'
$ws.Range("E48").Value = '  +10.57%  '

# Row 49
Set-TextValue $ws.Range("D49") '106.78'
$ws.Range("E49").Value = '  +2.42%  '

# Row 50
Set-TextValue $ws.Range("D50") '2.822.61'
$ws.Range("E50").Value = '  +0.48%  '

# Row 51
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D51") '1.68'
$ws.Range("E51").Value = '  +2.32%  '
